$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$techReq = "- Git: All code is in Git.`n- Go: Services are written in Go.`n- SQL: Extensive use of SQL with both Postgres and duckDB.`n- AWS: Use of AWS, but architecture remains cloud agnostic.`n- Svelte: Frontends are written in Svelte with JS and SCSS.`n- JavaScript: Frontends are written in Svelte with JS and SCSS.`n- SCSS: Frontends are written in Svelte with JS and SCSS.`n- Swift: Sprinkle of Swift for native where low level device access is needed.`n- Java: Sprinkle of Java for native where low level device access is needed."

$ws.Cells.Item(13, 1).Value = "Full Stack Engineer"
$ws.Cells.Item(13, 2).Value = $techReq
$ws.Cells.Item(13, 3).Value = 'No specific degree or years of experience are explicitly required, but a "can-do attitude and the curiosity to ask questions" are emphasized.'
$ws.Cells.Item(13, 4).Value = "https://www.linkedin.com/jobs/collections/recommended/?currentJobId=4168233573"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "4168233573"
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(13, 6).Value = "Payd"

$ws.Rows.Item(13).AutoFit()
